# Edit script: insert a new weekly price-group (Extra (muy buena)/Primera/Segunda,
# 2022-03-02, Region Metropolitana) at the top of the data table and shift the
# existing rows down by 3 positions (the 3 rows that fall past the old end of
# the table become new rows 106-108). Applies to worksheet "Sheet1" of the
# "Hortaliza, Vega Monumental Concepcion - Sandia" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared by every data row in this table.
$constA = 11
$constB = 'Vega Monumental Concepción'
$constC = 'Bíobío'
$constE = 8
$constF = 100112028
$constG = 'Sandia'
$constH = 'Sin especificar'
$constQ = 1
$constR = 'Hortaliza'

$rows = @(
  @{r=21; d=44622; i='Extra (muy buena)'; j=2000; k=3000; l=3000; m=3000; n='$/unidad'; o='Región Metropolitana'; p=3000},
  @{r=22; d=44622; i='Primera'; j=2500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región Metropolitana'; p=2500},
  @{r=23; d=44622; i='Segunda'; j=1000; k=2000; l=2000; m=2000; n='$/unidad'; o='Región Metropolitana'; p=2000},
  @{r=24; d=44211; i='Extra'; j=500; k=3500; l=3500; m=3500; n='$/unidad'; o='Región de O''Higgins'; p=3500},
  @{r=25; d=44211; i='Primera'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=26; d=44211; i='Segunda'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=27; d=44587; i='Primera'; j=1800; k=2000; l=2300; m=2133; n='$/unidad'; o='Región de O''Higgins'; p=2133},
  @{r=28; d=44587; i='Segunda'; j=2000; k=1400; l=1500; m=1450; n='$/unidad'; o='Región de O''Higgins'; p=1450},
  @{r=29; d=44594; i='Extra'; j=500; k=2800; l=2800; m=2800; n='$/unidad'; o='Región de O''Higgins'; p=2800},
  @{r=30; d=44594; i='Primera'; j=500; k=2400; l=2400; m=2400; n='$/unidad'; o='Región de O''Higgins'; p=2400},
  @{r=31; d=44594; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=32; d=44260; i='Extra'; j=300; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=33; d=44260; i='Primera'; j=400; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=34; d=44260; i='Segunda'; j=400; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=35; d=44574; i='Extra'; j=400; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=36; d=44574; i='Primera'; j=400; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=37; d=44574; i='Segunda'; j=400; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=38; d=44601; i='Primera'; j=800; k=2000; l=2500; m=2188; n='$/unidad'; o='Región Metropolitana'; p=2188},
  @{r=39; d=44601; i='Segunda'; j=700; k=1500; l=2000; m=1786; n='$/unidad'; o='Región Metropolitana'; p=1786},
  @{r=40; d=44546; i='Primera'; j=450; k=2500; l=3000; m=2778; n='$/unidad'; o='Región de O''Higgins'; p=2778},
  @{r=41; d=44202; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=42; d=44202; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=43; d=44202; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=44; d=44204; i='Extra'; j=500; k=3500; l=3500; m=3500; n='$/unidad'; o='Región de O''Higgins'; p=3500},
  @{r=45; d=44204; i='Primera'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=46; d=44204; i='Segunda'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=47; d=44244; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=48; d=44244; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=49; d=44244; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=50; d=44194; i='Extra'; j=400; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=51; d=44194; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=52; d=44194; i='Segunda'; j=400; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=53; d=44217; i='Extra'; j=500; k=3200; l=3200; m=3200; n='$/unidad'; o='Región de O''Higgins'; p=3200},
  @{r=54; d=44217; i='Primera'; j=500; k=2600; l=2600; m=2600; n='$/unidad'; o='Región de O''Higgins'; p=2600},
  @{r=55; d=44217; i='Segunda'; j=500; k=2200; l=2200; m=2200; n='$/unidad'; o='Región de O''Higgins'; p=2200},
  @{r=56; d=44596; i='Extra'; j=500; k=2800; l=2800; m=2800; n='$/unidad'; o='Región de O''Higgins'; p=2800},
  @{r=57; d=44596; i='Primera'; j=800; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=58; d=44596; i='Segunda'; j=600; k=2100; l=2100; m=2100; n='$/unidad'; o='Región de O''Higgins'; p=2100},
  @{r=59; d=44238; i='Extra'; j=400; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=60; d=44238; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=61; d=44238; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=62; d=44579; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=63; d=44579; i='Primera'; j=800; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=64; d=44579; i='Segunda'; j=800; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=65; d=44566; i='Primera'; j=1800; k=2200; l=2500; m=2367; n='$/unidad'; o='Paine'; p=2367},
  @{r=66; d=44566; i='Segunda'; j=1800; k=1800; l=2000; m=1933; n='$/unidad'; o='Paine'; p=1933},
  @{r=67; d=44196; i='Extra'; j=400; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=68; d=44196; i='Primera'; j=400; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=69; d=44196; i='Segunda'; j=400; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=70; d=44609; i='Extra'; j=500; k=2800; l=2800; m=2800; n='$/unidad'; o='Región de O''Higgins'; p=2800},
  @{r=71; d=44609; i='Primera'; j=500; k=2400; l=2400; m=2400; n='$/unidad'; o='Región de O''Higgins'; p=2400},
  @{r=72; d=44609; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=73; d=44225; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=74; d=44225; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=75; d=44225; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=76; d=44281; i='Primera'; j=400; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=77; d=44281; i='Segunda'; j=400; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=78; d=44568; i='Extra'; j=400; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=79; d=44568; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=80; d=44568; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=81; d=44511; i='Primera'; j=600; k=800; l=900; m=850; n='$/kilo (volumen en unidades)'; o='Perú'; p=850},
  @{r=82; d=44231; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=83; d=44231; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=84; d=44231; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=85; d=44573; i='Primera'; j=2700; k=2000; l=2200; m=2089; n='$/unidad'; o='Región de O''Higgins'; p=2089},
  @{r=86; d=44581; i='Extra'; j=400; k=2500; l=2500; m=2500; n='$/unidad'; o='Región del Maule'; p=2500},
  @{r=87; d=44581; i='Primera'; j=400; k=2000; l=2000; m=2000; n='$/unidad'; o='Región del Maule'; p=2000},
  @{r=88; d=44581; i='Segunda'; j=400; k=1500; l=1500; m=1500; n='$/unidad'; o='Región del Maule'; p=1500},
  @{r=89; d=44553; i='Extra'; j=500; k=3400; l=3400; m=3400; n='$/unidad'; o='Región de O''Higgins'; p=3400},
  @{r=90; d=44553; i='Primera'; j=500; k=2800; l=2800; m=2800; n='$/unidad'; o='Región de O''Higgins'; p=2800},
  @{r=91; d=44553; i='Segunda'; j=500; k=2400; l=2400; m=2400; n='$/unidad'; o='Región de O''Higgins'; p=2400},
  @{r=92; d=44208; i='Extra'; j=500; k=3500; l=3500; m=3500; n='$/kilo (volumen en unidades)'; o='Región de O''Higgins'; p=3500},
  @{r=93; d=44208; i='Primera'; j=500; k=3000; l=3000; m=3000; n='$/kilo (volumen en unidades)'; o='Región de O''Higgins'; p=3000},
  @{r=94; d=44208; i='Segunda'; j=500; k=2500; l=2500; m=2500; n='$/kilo (volumen en unidades)'; o='Región de O''Higgins'; p=2500},
  @{r=95; d=44617; i='Extra'; j=800; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=96; d=44617; i='Primera'; j=1000; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=97; d=44617; i='Segunda'; j=1000; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=98; d=44264; i='Extra'; j=300; k=2800; l=2800; m=2800; n='$/unidad'; o='Región de O''Higgins'; p=2800},
  @{r=99; d=44264; i='Primera'; j=300; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=100; d=44264; i='Segunda'; j=300; k=2200; l=2200; m=2200; n='$/unidad'; o='Región de O''Higgins'; p=2200},
  @{r=101; d=44232; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=102; d=44232; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=103; d=44232; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=104; d=44236; i='Extra'; j=500; k=3000; l=3000; m=3000; n='$/unidad'; o='Región de O''Higgins'; p=3000},
  @{r=105; d=44236; i='Primera'; j=500; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=106; d=44236; i='Segunda'; j=500; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000},
  @{r=107; d=44272; i='Primera'; j=300; k=2500; l=2500; m=2500; n='$/unidad'; o='Región de O''Higgins'; p=2500},
  @{r=108; d=44272; i='Segunda'; j=300; k=2000; l=2000; m=2000; n='$/unidad'; o='Región de O''Higgins'; p=2000}
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 1).Value = $constA
  $ws.Cells.Item($r, 2).Value = $constB
  $ws.Cells.Item($r, 3).Value = $constC
  $ws.Cells.Item($r, 4).Value = $row.d
  $ws.Cells.Item($r, 5).Value = $constE
  $ws.Cells.Item($r, 6).Value = $constF
  $ws.Cells.Item($r, 7).Value = $constG
  $ws.Cells.Item($r, 8).Value = $constH
  $ws.Cells.Item($r, 9).Value = $row.i
  $ws.Cells.Item($r, 10).Value = $row.j
  $ws.Cells.Item($r, 11).Value = $row.k
  $ws.Cells.Item($r, 12).Value = $row.l
  $ws.Cells.Item($r, 13).Value = $row.m
  $ws.Cells.Item($r, 14).Value = $row.n
  $ws.Cells.Item($r, 15).Value = $row.o
  $ws.Cells.Item($r, 16).Value = $row.p
  $ws.Cells.Item($r, 17).Value = $constQ
  $ws.Cells.Item($r, 18).Value = $constR
}
